$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Bump every due-date by one day (Wed -> Thu schedule shift)
$cells = @("B4", "B6", "B9", "B11", "B14", "B16", "B19", "B21", "B24", "B26", "B29", "B31", "B34", "B36")
foreach ($addr in $cells) {
    $ws.Range($addr).Value = $ws.Range($addr).Value2 + 1
}

# Move the active selection from F24 to B38
$ws.Range("B38").Select() | Out-Null
